# Updates the cryptocurrency price/volume table with the latest scraped
# values (GitHub Actions data refresh). For rows whose Price column holds a
# value that Excel would otherwise auto-parse as a number (losing formatting
# such as trailing zeros), the cell is forced to Text format before the
# value is written, then the style is reset back to Normal so no stray
# cell-style reference is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.302.24"
$ws.Range("E2").Value = "  -6.61%  "

# Row 3
$ws.Range("D3").Value = "3.493.14"
$ws.Range("E3").Value = "  -3.01%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "390.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.58%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "121.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.29%  "

# Row 7
$ws.Range("D7").Value = "3.485.05"
$ws.Range("E7").Value = "  -3.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -9.88%  "

# Row 9
$ws.Range("E9").Value = "  +0.05%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.675"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -12.21%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.152"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -14.80%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000332"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.47%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -9.68%  "

# Row 14
$ws.Range("D14").Value = "4.034.07"
$ws.Range("E14").Value = "  -2.86%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.17%  "

# Row 16
$ws.Range("E16").Value = "  -3.34%  "

# Row 17
$ws.Range("D17").Value = "3.491.42"
$ws.Range("E17").Value = "  -3.17%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -8.65%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.75%  "

# Row 20
$ws.Range("D20").Value = "63.331.08"
$ws.Range("E20").Value = "  -6.46%  "

# Row 21
$ws.Range("E21").Value = "  -11.79%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "394.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -15.17%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.96%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.76%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.79%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "33.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.68%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.30%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -11.90%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -14.78%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.24%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.49%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.110"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.90%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.91%  "

# Row 34
$ws.Range("E34").Value = "  -9.09%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.10%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "36.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -11.78%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.10%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0435"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -12.30%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.994"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.50%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.72"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +15.88%  "

# Row 41
$ws.Range("D41").Value = "0.0₃0637"
$ws.Range("E41").Value = "  -11.63%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.131"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.63%  "

# Row 43
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +12.11%  "

# Row 44
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "140.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.30%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.93%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "24.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +14.69%  "

# Row 47
$ws.Range("B47").Value = "LidoDAOToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.52%  "

# Row 48
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.47%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.43"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -11.06%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.83%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.274"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -12.04%  "
